# Auto-generated edit script for cryptos.xlsx update
# Commit: Updated cryptos list on Tue Jun  4 14:42:20 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextCell 2 4 '69.410.00'  # D2: '69.760.68' -> '69.410.00'
Set-TextCell 2 5 '  +0.18%  '  # E2: '  +0.10%  ' -> '  +0.18%  '
Set-TextCell 3 4 '3.774.23'  # D3: '3.784.92' -> '3.774.23'
Set-TextCell 3 5 '  -0.32%  '  # E3: '  -0.49%  ' -> '  -0.32%  '
Set-TextCell 4 4 '1.00'  # D4: '0.999' -> '1.00'
Set-TextCell 4 5 '  +0.08%  '  # E4: '  -0.25%  ' -> '  +0.08%  '
Set-TextCell 5 4 '662.13'  # D5: '662.30' -> '662.13'
Set-TextCell 5 5 '  +5.39%  '  # E5: '  +4.86%  ' -> '  +5.39%  '
Set-TextCell 6 4 '166.12'  # D6: '166.94' -> '166.12'
Set-TextCell 6 5 '  +1.34%  '  # E6: '  +1.37%  ' -> '  +1.34%  '
Set-TextCell 7 4 '3.774.18'  # D7: '3.782.89' -> '3.774.18'
Set-TextCell 7 5 '  -0.19%  '  # E7: '  -0.25%  ' -> '  -0.19%  '
Set-TextCell 8 5 '  +0.02%  '  # E8: '  -0.06%  ' -> '  +0.02%  '
Set-TextCell 9 4 '0.526'  # D9: '0.527' -> '0.526'
Set-TextCell 9 5 '  +1.33%  '  # E9: '  +1.31%  ' -> '  +1.33%  '
Set-TextCell 10 4 '0.159'  # D10: '0.160' -> '0.159'
Set-TextCell 10 5 '  -0.74%  '  # E10: '  -0.71%  ' -> '  -0.74%  '
Set-TextCell 11 4 '0.457'  # D11: '0.458' -> '0.457'
Set-TextCell 11 5 '  +1.04%  '  # E11: '  +1.08%  ' -> '  +1.04%  '
Set-TextCell 12 4 '6.93'  # D12: '6.94' -> '6.93'
Set-TextCell 12 5 '  +4.44%  '  # E12: '  +4.27%  ' -> '  +4.44%  '
Set-TextCell 13 5 '  -2.82%  '  # E13: '  -3.14%  ' -> '  -2.82%  '
Set-TextCell 14 4 '35.04'  # D14: '35.21' -> '35.04'
Set-TextCell 14 5 '  -1.40%  '  # E14: '  -1.32%  ' -> '  -1.40%  '
Set-TextCell 15 4 '4.421.31'  # D15: '4.425.95' -> '4.421.31'
Set-TextCell 15 5 '  -0.46%  '  # E15: '  -0.68%  ' -> '  -0.46%  '
Set-TextCell 16 4 '3.785.45'  # D16: '3.786.80' -> '3.785.45'
Set-TextCell 16 5 '  +0.49%  '  # E16: '  -1.77%  ' -> '  +0.49%  '
Set-TextCell 17 4 '69.487.00'  # D17: '69.693.39' -> '69.487.00'
Set-TextCell 17 5 '  +0.17%  '  # E17: '  -0.06%  ' -> '  +0.17%  '
Set-TextCell 18 4 '17.68'  # D18: '17.74' -> '17.68'
Set-TextCell 18 5 '  -1.53%  '  # E18: '  -1.39%  ' -> '  -1.53%  '
Set-TextCell 19 2 'TRON'  # B19: 'Polkadot' -> 'TRON'
Set-TextCell 19 3 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'  # C19: 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' -> 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell 19 4 '0.114'  # D19: '7.13' -> '0.114'
Set-TextCell 19 5 '  +0.45%  '  # E19: '  +0.07%  ' -> '  +0.45%  '
Set-TextCell 20 2 'Polkadot'  # B20: 'TRON' -> 'Polkadot'
Set-TextCell 20 3 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'  # C20: 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx' -> 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 20 4 '7.08'  # D20: '0.114' -> '7.08'
Set-TextCell 20 5 '  +0.01%  '  # E20: '  +0.31%  ' -> '  +0.01%  '
Set-TextCell 21 4 '469.20'  # D21: '470.86' -> '469.20'
Set-TextCell 21 5 '  +0.23%  '  # E21: '  +0.57%  ' -> '  +0.23%  '
Set-TextCell 22 4 '9.60'  # D22: '9.64' -> '9.60'
Set-TextCell 22 5 '  -0.09%  '  # E22: '  +0.10%  ' -> '  -0.09%  '
Set-TextCell 23 4 '0.709'  # D23: '0.712' -> '0.709'
Set-TextCell 23 5 '  +0.86%  '  # E23: '  +1.01%  ' -> '  +0.86%  '
Set-TextCell 24 4 '0.0000144'  # D24: '0.0000145' -> '0.0000144'
Set-TextCell 24 5 '  -3.07%  '  # E24: '  -3.02%  ' -> '  -3.07%  '
Set-TextCell 25 4 '82.31'  # D25: '82.15' -> '82.31'
Set-TextCell 25 5 '  -1.15%  '  # E25: '  -1.80%  ' -> '  -1.15%  '
Set-TextCell 26 4 '12.25'  # D26: '12.30' -> '12.25'
Set-TextCell 26 5 '  +1.79%  '  # E26: '  +1.76%  ' -> '  +1.79%  '
Set-TextCell 27 4 '10.24'  # D27: '10.32' -> '10.24'
Set-TextCell 27 5 '  +2.35%  '  # E27: '  +2.65%  ' -> '  +2.35%  '
Set-TextCell 28 4 '2.12'  # D28: '2.13' -> '2.12'
Set-TextCell 28 5 '  -1.37%  '  # E28: '  -1.64%  ' -> '  -1.37%  '
Set-TextCell 29 5 '  +0.05%  '  # E29: '  +0.09%  ' -> '  +0.05%  '
Set-TextCell 30 4 '3.931.76'  # D30: '3.933.60' -> '3.931.76'
Set-TextCell 30 5 '  -0.60%  '  # E30: '  -0.83%  ' -> '  -0.60%  '
Set-TextCell 31 4 '2.78'  # D31: '2.79' -> '2.78'
Set-TextCell 31 5 '  +3.98%  '  # E31: '  +3.77%  ' -> '  +3.98%  '
Set-TextCell 32 4 '2.27'  # D32: '2.28' -> '2.27'
Set-TextCell 32 5 '  +2.89%  '  # E32: '  +2.98%  ' -> '  +2.89%  '
Set-TextCell 33 4 '7.21'  # D33: '7.26' -> '7.21'
Set-TextCell 33 5 '  -0.79%  '  # E33: '  -0.66%  ' -> '  -0.79%  '
Set-TextCell 34 4 '28.81'  # D34: '28.88' -> '28.81'
Set-TextCell 34 5 '  -0.44%  '  # E34: '  -0.71%  ' -> '  -0.44%  '
Set-TextCell 35 5 '  +16.70%  '  # E35: '  +15.88%  ' -> '  +16.70%  '
Set-TextCell 36 5 '  +0.13%  '  # E36: '  +0.05%  ' -> '  +0.13%  '
Set-TextCell 37 4 '3.737.97'  # D37: '3.740.49' -> '3.737.97'
Set-TextCell 37 5 '  +0.01%  '  # E37: '  -0.35%  ' -> '  +0.01%  '
Set-TextCell 38 4 '8.89'  # D38: '8.94' -> '8.89'
Set-TextCell 38 5 '  -1.00%  '  # E38: '  -0.95%  ' -> '  -1.00%  '
Set-TextCell 39 5 '  -1.61%  '  # E39: '  -1.43%  ' -> '  -1.61%  '
Set-TextCell 40 4 '3.29'  # D40: '3.32' -> '3.29'
Set-TextCell 40 5 '  -0.50%  '  # E40: '  -0.08%  ' -> '  -0.50%  '
Set-TextCell 41 4 '5.84'  # D41: '5.87' -> '5.84'
Set-TextCell 41 5 '  -0.29%  '  # E41: '  -0.06%  ' -> '  -0.29%  '
Set-TextCell 42 5 '  +0.15%  '  # E42: '  +0.02%  ' -> '  +0.15%  '
Set-TextCell 43 4 '0.959'  # D43: '0.961' -> '0.959'
Set-TextCell 43 5 '  -0.91%  '  # E43: '  -1.49%  ' -> '  -0.91%  '
Set-TextCell 44 5 '  -0.04%  '  # E44: '  -0.03%  ' -> '  -0.04%  '
Set-TextCell 45 4 '45.81'  # D45: '46.22' -> '45.81'
Set-TextCell 45 5 '  +7.63%  '  # E45: '  +7.22%  ' -> '  +7.63%  '
Set-TextCell 46 4 '2.02'  # D46: '2.04' -> '2.02'
Set-TextCell 46 5 '  +4.49%  '  # E46: '  +4.32%  ' -> '  +4.49%  '
Set-TextCell 47 4 '157.69'  # D47: '157.87' -> '157.69'
Set-TextCell 47 5 '  +1.69%  '  # E47: '  +2.03%  ' -> '  +1.69%  '
Set-TextCell 48 4 '47.92'  # D48: '48.19' -> '47.92'
Set-TextCell 48 5 '  +2.02%  '  # E48: '  +2.43%  ' -> '  +2.02%  '
Set-TextCell 49 4 '0.299'  # D49: '0.301' -> '0.299'
Set-TextCell 49 5 '  +0.50%  '  # E49: '  +0.71%  ' -> '  +0.50%  '
Set-TextCell 50 4 '1.40'  # D50: '1.41' -> '1.40'
Set-TextCell 50 5 '  +0.77%  '  # E50: '  +1.16%  ' -> '  +0.77%  '
Set-TextCell 51 4 '8.44'  # D51: '8.47' -> '8.44'
Set-TextCell 51 5 '  +0.34%  '  # E51: '  +0.57%  ' -> '  +0.34%  '
